$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (ID_Location) values for specific rows
$ws.Range("B2").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2
$ws.Range("B14").Value = 1
$ws.Range("B17").Value = 2

# Update the active selection to F16
$ws.Range("F16").Select()
